$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")

# Row 76
$ws.Range("H76").Value = 4133.6665
$ws.Range("J76").Value = 5000
$ws.Range("L76").Value = 5000
$ws.Range("N76").Value = -5630

# Row 79
$ws.Range("H79").Value = 4133.6665
$ws.Range("J79").Value = 5000
$ws.Range("L79").Value = 5000
$ws.Range("N79").Value = -7184

# Row 98
$ws.Range("H98").Value = 1637.75
$ws.Range("I98").Value = 1637.75
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 1637.75
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()

# Row 99
$ws.Range("H99").Value = 1919.125
$ws.Range("I99").Value = 654.6
$ws.Range("K99").Value = 1963.8
$ws.Range("M99").Value = -465.8000000000002

# Row 122
$ws.Range("H122").Value = 1637.75
$ws.Range("I122").Value = 1637.75
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4913.25
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

# Row 129
$ws.Range("H129").Value = 1453.5
$ws.Range("I129").Value = 735.6667
$ws.Range("K129").Value = 2207.0001
$ws.Range("M129").Value = 2792.9999

# Row 132
$ws.Range("H132").Value = 22016.1
$ws.Range("I132").Value = 26920.125
$ws.Range("K132").Value = 80760.375
$ws.Range("M132").Value = -78230.375

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")

# Row 61
$ws.Range("H61").Value = 7499.5
$ws.Range("I61").Value = 8500
$ws.Range("K61").Value = 8500
$ws.Range("M61").Value = -8288

# Row 74
$ws.Range("H74").Value = 6929.846
$ws.Range("I74").Value = 3565
$ws.Range("K74").Value = 3565
$ws.Range("M74").Value = -2691

# Row 77
$ws.Range("H77").Value = 6929.846
$ws.Range("I77").Value = 3565
$ws.Range("K77").Value = 17825
$ws.Range("M77").Value = -13457

# Row 136
$ws.Range("H136").Value = 7499.5
$ws.Range("I136").Value = 8500
$ws.Range("K136").Value = 25500
$ws.Range("M136").Value = -22950

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")

# Row 7
$ws.Range("H7").Value = 7341729.5
$ws.Range("I7").Value = 6666745.5
$ws.Range("J7").Value = 10041667
$ws.Range("K7").Value = 6666745.5
$ws.Range("L7").Value = 10041667
$ws.Range("M7").Value = -6666632.5
$ws.Range("N7").Value = -10041893

# Row 105
$ws.Range("H105").Value = 1866.8
$ws.Range("J105").Value = 2196.4
$ws.Range("L105").Value = 2196.4
$ws.Range("N105").Value = -5690.4

# Row 107
$ws.Range("H107").Value = 4802.75
$ws.Range("I107").Value = 1511.6364
$ws.Range("K107").Value = 1511.6364
$ws.Range("M107").Value = 408.3635999999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")

# Row 58
$ws.Range("H58").Value = 2329.923
$ws.Range("I58").Value = 1030
$ws.Range("K58").Value = 1030
$ws.Range("M58").Value = -827

# Row 88
$ws.Range("H88").Value = 10208.2
$ws.Range("J88").Value = 11125
$ws.Range("L88").Value = 11125
$ws.Range("N88").Value = -11937

# Row 91
$ws.Range("H91").Value = 10208.2
$ws.Range("J91").Value = 11125
$ws.Range("L91").Value = 11125
$ws.Range("N91").Value = -13933

# Row 134
$ws.Range("H134").Value = 4664.8335
$ws.Range("J134").Value = 8000
$ws.Range("L134").Value = 24000
$ws.Range("N134").Value = -29070

# Row 136
$ws.Range("H136").Value = 2329.923
$ws.Range("I136").Value = 1030
$ws.Range("K136").Value = 3090
$ws.Range("M136").Value = -540

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")

# Row 97
$ws.Range("H97").Value = 7983
$ws.Range("I97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("M97").ClearContents()

# Row 98
$ws.Range("H98").Value = 248.33333
$ws.Range("I98").Value = 122.5
$ws.Range("K98").Value = 367.5
$ws.Range("M98").Value = 1130.5

# Row 122
$ws.Range("H122").Value = 705.5
$ws.Range("I122").Value = 749
$ws.Range("J122").Value = 679.4
$ws.Range("K122").Value = 6741
$ws.Range("L122").Value = 6114.599999999999
$ws.Range("M122").Value = -4291
$ws.Range("N122").Value = -11014.6

# Row 123
$ws.Range("H123").Value = 472.5
$ws.Range("I123").Value = 472.5
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 1417.5
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

# Row 131
$ws.Range("H131").Value = 2038.7333
$ws.Range("J131").Value = 2819.8
$ws.Range("L131").Value = 8459.400000000001
$ws.Range("N131").Value = -18539.4

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")

# Row 70
$ws.Range("H70").Value = 5185.8
$ws.Range("I70").Value = 4708.1816
$ws.Range("K70").Value = 4708.1816
$ws.Range("M70").Value = -4438.1816

# Row 73
$ws.Range("H73").Value = 5185.8
$ws.Range("I73").Value = 4708.1816
$ws.Range("K73").Value = 4708.1816
$ws.Range("M73").Value = -3772.1816

# Row 80
$ws.Range("H80").Value = 2350
$ws.Range("I80").Value = 1966.6666
$ws.Range("J80").Value = 3500
$ws.Range("K80").Value = 1966.6666
$ws.Range("L80").Value = 3500
$ws.Range("M80").Value = -968.6666
$ws.Range("N80").Value = -5496

# Row 83
$ws.Range("H83").Value = 2350
$ws.Range("I83").Value = 1966.6666
$ws.Range("J83").Value = 3500
$ws.Range("K83").Value = 9833.333000000001
$ws.Range("L83").Value = 17500
$ws.Range("M83").Value = -4841.333000000001
$ws.Range("N83").Value = -27484

# Row 140
$ws.Range("H140").Value = 129949.5
$ws.Range("J140").Value = 100000
$ws.Range("L140").Value = 100000
$ws.Range("N140").Value = -110360

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")

# Row 22
$ws.Range("H22").Value = 1625
$ws.Range("J22").Value = 1937.5
$ws.Range("L22").Value = 1937.5
$ws.Range("N22").Value = -2527.5

# Row 27
$ws.Range("H27").Value = 1625
$ws.Range("J27").Value = 1937.5
$ws.Range("L27").Value = 1937.5
$ws.Range("N27").Value = -2151.5

# Row 53
$ws.Range("H53").Value = 19990
$ws.Range("I53").Value = 19990
$ws.Range("K53").Value = 19990
$ws.Range("M53").Value = -19472

# Row 74
$ws.Range("H74").Value = 33995
$ws.Range("I74").Value = 33995
$ws.Range("K74").Value = 33995
$ws.Range("M74").Value = -32997

# Row 77
$ws.Range("H77").Value = 33995
$ws.Range("I77").Value = 33995
$ws.Range("K77").Value = 101985
$ws.Range("M77").Value = -96993

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")

# Row 2
$ws.Range("H2").Value = 157714.77
$ws.Range("I2").Value = 170024.08
$ws.Range("K2").Value = 170024.08
$ws.Range("M2").Value = -169912.08

# Row 4
$ws.Range("H4").Value = 176797.92
$ws.Range("I4").Value = 212077.5
$ws.Range("J4").Value = 400
$ws.Range("K4").Value = 212077.5
$ws.Range("L4").Value = 400
$ws.Range("M4").Value = -211964.5
$ws.Range("N4").Value = -626
